$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.413.33"
$ws.Range("E2").Value = "  -0.62%  "

$ws.Range("D3").Value = "1.819.27"
$ws.Range("E3").Value = "  -2.22%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -1.30%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "332.35"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.80%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.002"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.12%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4564"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.83%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3808"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.71%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "45.93"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.01%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07838"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.47%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.9540"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -4.49%  "

$ws.Range("E12").Value = "  -3.21%  "

$ws.Range("D13").Value = "1.837.55"
$ws.Range("E13").Value = "  -1.38%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.810"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.16%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.046"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.16%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.004"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.21%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "89.14"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.09%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06565"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.28%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.00001016"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.49%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.04"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.82%  "

$ws.Range("E21").Value = "  -1.02%  "

$ws.Range("D22").Value = "27.397.39"
$ws.Range("E22").Value = "  -0.76%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.274"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.20%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "10.78"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.30%  "

$ws.Range("E25").Value = "  -2.21%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "158.11"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.93%  "

$ws.Range("D27").Value = "2.045.50"
$ws.Range("E27").Value = "  -1.97%  "

$ws.Range("E28").Value = "  -1.77%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.034"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -4.82%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.248"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.21%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "117.59"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.38%  "

$ws.Range("E32").Value = "  -1.43%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.9285"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -4.64%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.562"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.76%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.195"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.96%  "

$ws.Range("E36").Value = "  -1.74%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.05888"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.22%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02175"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.48%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "8.061"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.97%  "

$ws.Range("E40").Value = "  -1.13%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.135"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -5.18%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.5726"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.42%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1809"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.06%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "9.871"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -4.00%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.263"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.85%  "

$ws.Range("E46").Value = "  -3.67%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "11.73"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.65%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.868"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.70%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "110.04"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.48%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06560"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.61%  "

$ws.Range("E51").Value = "  -33.30%  "
